$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1. Remove the old "Moritz Passenbrunner" (E:G) and "Paul Achleitener" (I:K)
#    tables entirely - they get rebuilt, shifted right, with new data below.
# ---------------------------------------------------------------------------
$ws.Range("E1:K8").Clear()

# ---------------------------------------------------------------------------
# 2. Table 1 "Jannis Katsanis" (A:D) - add the new "Issue Number" column (D)
#    to the existing header band / column-header row. A1:C1 and A3:C3 stay.
# ---------------------------------------------------------------------------
$ws.Range("A1").Copy()
$ws.Range("D1").PasteSpecial(-4122)

$ws.Range("A3").Copy()
$ws.Range("D3").PasteSpecial(-4122)
$ws.Range("D3").Value = "Issue Number"

# ---------------------------------------------------------------------------
# 3. Table 2 "Moritz Passenbrunner" now lives in F:I (was E:G).
# ---------------------------------------------------------------------------
$ws.Range("A1").Copy()
$ws.Range("F1:I1").PasteSpecial(-4122)
$ws.Range("F1").Value = "Moritz Passenbrunner"

$ws.Range("A3").Copy()
$ws.Range("F3:I3").PasteSpecial(-4122)
$ws.Range("F3").Value = "Tätigkeit"
$ws.Range("G3").Value = "Zeit"
$ws.Range("H3").Value = "Datum"
$ws.Range("I3").Value = "Issue Number"

$ws.Range("F4").Value = "Projektauftrag"
$ws.Range("G4").Value = "1h"
$ws.Range("H4").NumberFormat = "d-mmm"
$ws.Range("H4").Value = 45364

$ws.Range("F5").Value = "Mockup"
$ws.Range("G5").Value = "2h"
$ws.Range("H5").NumberFormat = "d-mmm"
$ws.Range("H5").Value = 45371

$ws.Range("F6").Value = "Server Setup"
$ws.Range("G6").Value = "3h"
$ws.Range("H6").NumberFormat = "d-mmm"
$ws.Range("H6").Value = 45385

$ws.Range("F7").Value = "Server Setup"
$ws.Range("G7").Value = "3h"
$ws.Range("H7").NumberFormat = "d-mmm"
$ws.Range("H7").Value = 45392

$ws.Range("F8").Value = "Server Setup"
$ws.Range("G8").Value = "2h"
$ws.Range("H8").NumberFormat = "d-mmm"
$ws.Range("H8").Value = 45399

# ---------------------------------------------------------------------------
# 4. New separator cell J3 - highlighted (white) header-row filler between
#    table 2 and table 3.
# ---------------------------------------------------------------------------
$ws.Range("J3").Interior.ThemeColor = 2

# ---------------------------------------------------------------------------
# 5. Table 3 "Paul Achleitener" now lives in K:N (was I:K), gained a new
#    "Issue Number" column (N) and a new row (9).
# ---------------------------------------------------------------------------
$ws.Range("A1").Copy()
$ws.Range("K1:N1").PasteSpecial(-4122)
$ws.Range("K1").Value = "Paul Achleitener"

$ws.Range("A3").Copy()
$ws.Range("K3:N3").PasteSpecial(-4122)
$ws.Range("K3").Value = "Tätigkeit"
$ws.Range("L3").Value = "Zeit"
$ws.Range("M3").Value = "Datum"
$ws.Range("N3").Value = "Issue Number"

$ws.Range("K4").Value = "Mockups gemacht"
$ws.Range("L4").Value = "2h"
$ws.Range("M4").NumberFormat = "d-mmm"
$ws.Range("M4").Value = 45357

$ws.Range("K5").Value = "Projektauftrag finalisiert"
$ws.Range("L5").Value = "1h"
$ws.Range("M5").NumberFormat = "d-mmm"
$ws.Range("M5").Value = 45364

$ws.Range("K6").Value = "Pflichtenheft weiter geschrieben"
$ws.Range("L6").Value = "1h"
$ws.Range("M6").NumberFormat = "d-mmm"
$ws.Range("M6").Value = 45364

$ws.Range("K7").Value = "Statische Inhalte verbessert"
$ws.Range("L7").Value = "1h"
$ws.Range("M7").NumberFormat = "d-mmm"
$ws.Range("M7").Value = 45385
$ws.Range("N7").Value = 19

$ws.Range("K8").Value = "Projekt speichern angefangen"
$ws.Range("L8").Value = "3h"
$ws.Range("M8").NumberFormat = "d-mmm"
$ws.Range("M8").Value = 45399
$ws.Range("N8").Value = 26

$ws.Range("K9").Value = "Organisatorische Sachen geupdatet"
$ws.Range("L9").Value = "1h"
$ws.Range("M9").NumberFormat = "d-mmm"
$ws.Range("M9").Value = 45406

# ---------------------------------------------------------------------------
# 6. Column widths: the old column-I custom width (table 2 used to live
#    there) is no longer meaningful now that table 2 moved to F:I headers
#    only; give the new table-3 "Tätigkeit" column (K) and "Issue Number"
#    column (N) their own custom widths.
# ---------------------------------------------------------------------------
$ws.Columns.Item(9).ColumnWidth = $ws.StandardWidth
$ws.Columns.Item(11).ColumnWidth = 29.5
$ws.Columns.Item(14).ColumnWidth = 13

# ---------------------------------------------------------------------------
# 7. Selection, matching the author's final cursor position.
# ---------------------------------------------------------------------------
$ws.Range("M13").Select()
